$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 4: Mary's Math score is replaced by her English score
$ws.Range("C4").Value = "English"
$ws.Range("D4").Value = 67

# Prepare new rows 5-7 with the same look (style) as the existing data rows
$ws.Range("A4:D4").Copy()
$ws.Range("A5:D7").PasteSpecial(-4122)

# Row 5: Mary's Math score (previously on row 4, now moved down)
$ws.Range("A5").Value = 102
$ws.Range("B5").Value = "Mary"
$ws.Range("C5").Value = "Math"
$ws.Range("D5").Value = 35

# Row 6: new student Sidhaarth, Social Science
$ws.Range("A6").Value = 103
$ws.Range("B6").Value = "Sidhaarth"
$ws.Range("C6").Value = "Social Science"
$ws.Range("D6").Value = 76

# Row 7: John, Social Science
$ws.Range("A7").Value = 101
$ws.Range("B7").Value = "John"
$ws.Range("C7").Value = "Social Science"
$ws.Range("D7").Value = 68

# Rows 6 and 7 wrap onto two lines, so they need the taller row height
$ws.Rows.Item(6).RowHeight = 29
$ws.Rows.Item(7).RowHeight = 29

# Column A is widened slightly
$ws.Columns.Item(1).ColumnWidth = 10

# Move the active selection to D4, matching the latest edit location
$ws.Range("D4").Select() | Out-Null
